$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.6073882409096184, 0.9996992968541456, 0.8120463847722224, 0.5179604801884602, 0.745758354755784, 0.7210935702783527, 0.6685034612559557, 0.7375851701384918, 0.8760269644770509, 0.8483591299835292, 0.8438772889068291, 0.7275429459198581),
    @(0.6052851972438196, 0.9997328637169386, 0.8046568310459208, 0.5204117365252723, 0.7154426413881748, 0.7328819624166091, 0.6580486743326707, 0.7566404121797502, 0.882435386740702, 0.8565469106067102, 0.8474757807052953, 0.7359225986473735),
    @(0.6003187955832799, 0.9996526000960423, 0.8043230387084995, 0.55448396667161, 0.7410065874035989, 0.732277199074074, 0.6468706658286273, 0.729675193209258, 0.8760049833763268, 0.8525815217391304, 0.831045191002165, 0.7219258476592869),
    @(0.6134154699060034, 0.9996902201004312, 0.8102067230151242, 0.580200805801022, 0.7262194730077121, 0.7190752961812744, 0.6786450693445703, 0.7599295097294061, 0.8763596614761464, 0.8464603404706571, 0.8495364848996367, 0.7204239941382358),
    @(0.5915305691245603, 0.9996948155637898, 0.8092937867076474, 0.5170893019781647, 0.7057085475578406, 0.720903514636531, 0.6673467410940583, 0.7192313716295878, 0.8926487747957992, 0.8547440124320299, 0.8676046291879197, 0.7151544111117996),
    @(0.6012017235980986, 0.9996659583373581, 0.7944362237911202, 0.5537898193275778, 0.7088753213367609, 0.7070847171900161, 0.6761083642204693, 0.7318025849635491, 0.8799027168640494, 0.8433428143543693, 0.8663535283108876, 0.7145139993062471),
    @(0.6176899800533808, 0.9996201464200257, 0.7911933979489365, 0.5129803223980539, 0.7166002570694088, 0.7198206917126754, 0.6730947237836521, 0.7342374588142145, 0.8814447469640001, 0.8454711085038694, 0.8676011889471946, 0.7423874041059181),
    @(0.5983164703023085, 0.9996751207208244, 0.8052611758619639, 0.55814415747672, 0.7421850899742931, 0.7360419182769727, 0.6750923790688814, 0.7301056681441822, 0.8819450964262456, 0.8543410572866136, 0.849778448497303, 0.7242996833395191),
    @(0.6107201583468171, 0.9996549691858483, 0.8060614215939869, 0.5457155036334033, 0.6886278920308483, 0.7239767548021624, 0.6687675565839443, 0.7432534104088454, 0.8836398765481774, 0.8619407983800005, 0.8652480642912187, 0.7214215033994529),
    @(0.6220509766399034, 0.9996611916144953, 0.8026279252000852, 0.5181771100642824, 0.7436945694087405, 0.7207085166494134, 0.6545867668259781, 0.7236716052775458, 0.8710339015634152, 0.8531558121432279, 0.867450965102198, 0.7245152833239321),
    @(0.6351509352577843, 0.9997293528971057, 0.7983675557651928, 0.5443070212149752, 0.7219633676092545, 0.7229330213365539, 0.66639243610131, 0.7132489230920015, 0.8760534163101257, 0.8480221001331197, 0.8656517192029651, 0.736885931318469),
    @(0.6074506425176094, 0.9996730941500266, 0.8031038883478895, 0.5649132587577869, 0.7342079048843188, 0.7078130391649412, 0.6519280044861457, 0.7457718809544628, 0.8792425387202677, 0.8551388591187021, 0.8314947157902463, 0.7266664920572965),
    @(0.6235518995268425, 0.9997175930778278, 0.8015447643102767, 0.514497507979977, 0.7337355398457583, 0.7107106014205198, 0.6618315792857177, 0.7247103182343608, 0.8822315958577175, 0.8653463510525484, 0.8483427213680232, 0.7169635345365624),
    @(0.6103889919184443, 0.9996735508420374, 0.8069106425550044, 0.5444607585462038, 0.7406041131105399, 0.7004260301644811, 0.6595090172482917, 0.7386046654642836, 0.878503377638291, 0.8419887017440942, 0.8310967946130418, 0.7305555963687214),
    @(0.6013358323170264, 0.9996465489268994, 0.8008025280727916, 0.5701737154197767, 0.7306820372750642, 0.7257937406544743, 0.6743874009153856, 0.7364186599353711, 0.8856673536692045, 0.8457626604769748, 0.8447556970386408, 0.7138288704668906),
    @(0.6070340844151436, 0.9996480902624357, 0.8045864588722192, 0.5415537253738807, 0.7296392994858611, 0.7011435688405797, 0.6692579572834494, 0.73392997671784, 0.8773938908932672, 0.8515728999796935, 0.8492945213019706, 0.7233944828123693),
    @(0.598943770677375, 0.9997126836387118, 0.802415857703661, 0.5545336595463506, 0.7066749678663239, 0.7185837574764206, 0.6562582122570638, 0.7467212319270191, 0.8874847063697504, 0.8515972253559261, 0.8485995926754982, 0.7187554099625718),
    @(0.5988266308167605, 0.9997405418513702, 0.7946107277623932, 0.5510318647794762, 0.6732029241645244, 0.7254715896020243, 0.6783357998157414, 0.7258191755444114, 0.8803244559491289, 0.8569079132916676, 0.8498885362005065, 0.7436931414929996),
    @(0.6219404053696036, 0.9996784888244039, 0.8048898199993914, 0.5305048252557825, 0.7238303341902315, 0.7174977894237404, 0.650277408686217, 0.7137697208927358, 0.8741250404228716, 0.8558464102795514, 0.8832623114748083, 0.71912017319546),
    @(0.5928076125582689, 0.999637643432689, 0.8035812779586744, 0.5204451239254886, 0.7202699228791773, 0.7136342233149299, 0.6423719235391903, 0.7399364472942056, 0.8830743289058181, 0.8620204728007039, 0.8488966001247661, 0.7137004686983957),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}